# Auto-generated edit script: applies market-price + profit-column refresh
# to the Unicorn_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1169.5555
$ws.Range("I17").Value = 862.5
$ws.Range("J17").Value = 1199.5122
$ws.Range("K17").Value = 2587.5
$ws.Range("L17").Value = 3598.536599999999
$ws.Range("M17").Value = -2419.5
$ws.Range("N17").Value = -3934.536599999999
$ws.Range("H19").Value = 2916.85
$ws.Range("I19").Value = 3938.3572
$ws.Range("J19").Value = 533.3333
$ws.Range("K19").Value = 3938.3572
$ws.Range("L19").Value = 533.3333
$ws.Range("M19").Value = -3763.3572
$ws.Range("N19").Value = -883.3333
$ws.Range("H53").Value = 201.52942
$ws.Range("I53").Value = 146.42857
$ws.Range("J53").Value = 240.1
$ws.Range("K53").Value = 146.42857
$ws.Range("L53").Value = 240.1
$ws.Range("M53").Value = 490.57143
$ws.Range("N53").Value = -1514.1
$ws.Range("H86").Value = 5316.8184
$ws.Range("I86").Value = 2807.9
$ws.Range("J86").Value = 7407.5835
$ws.Range("K86").Value = 2807.9
$ws.Range("L86").Value = 7407.5835
$ws.Range("M86").Value = -1684.9
$ws.Range("N86").Value = -9653.583500000001
$ws.Range("H89").Value = 5316.8184
$ws.Range("I89").Value = 2807.9
$ws.Range("J89").Value = 7407.5835
$ws.Range("K89").Value = 14039.5
$ws.Range("L89").Value = 37037.9175
$ws.Range("M89").Value = -8423.5
$ws.Range("N89").Value = -48269.9175
$ws.Range("H135").Value = 955.1667
$ws.Range("I135").Value = 353
$ws.Range("J135").Value = 2520.8
$ws.Range("K135").Value = 3177
$ws.Range("L135").Value = 22687.2
$ws.Range("M135").Value = -642
$ws.Range("N135").Value = -27757.2
$ws.Range("H137").Value = 55722.85
$ws.Range("I137").Value = 83521.766
$ws.Range("J137").Value = 4096.2856
$ws.Range("K137").Value = 250565.298
$ws.Range("L137").Value = 12288.8568
$ws.Range("M137").Value = -248015.298
$ws.Range("N137").Value = -17388.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 550.36365
$ws.Range("I2").Value = 484.35715
$ws.Range("K2").Value = 484.35715
$ws.Range("M2").Value = -371.35715
$ws.Range("H45").Value = 1785.9286
$ws.Range("I45").Value = 1214.7142
$ws.Range("J45").Value = 2357.1428
$ws.Range("K45").Value = 1214.7142
$ws.Range("L45").Value = 2357.1428
$ws.Range("M45").Value = -837.7141999999999
$ws.Range("N45").Value = -3111.1428
$ws.Range("H53").Value = 3900
$ws.Range("I53").Value = 3900
$ws.Range("K53").Value = 3900
$ws.Range("M53").Value = -3218
$ws.Range("H61").Value = 3575.4736
$ws.Range("I61").Value = 2715.6365
$ws.Range("J61").Value = 4757.75
$ws.Range("K61").Value = 2715.6365
$ws.Range("L61").Value = 4757.75
$ws.Range("M61").Value = -2503.6365
$ws.Range("N61").Value = -5181.75
$ws.Range("H74").Value = 3391.3914
$ws.Range("I74").Value = 1840
$ws.Range("K74").Value = 1840
$ws.Range("M74").Value = -966
$ws.Range("H77").Value = 3391.3914
$ws.Range("I77").Value = 1840
$ws.Range("K77").Value = 9200
$ws.Range("M77").Value = -4832
$ws.Range("H116").Value = 550.36365
$ws.Range("I116").Value = 484.35715
$ws.Range("K116").Value = 484.35715
$ws.Range("M116").Value = 1809.64285
$ws.Range("H122").Value = 2045
$ws.Range("I122").Value = 2404.9546
$ws.Range("J122").Value = 1517.0667
$ws.Range("K122").Value = 7214.8638
$ws.Range("L122").Value = 4551.2001
$ws.Range("M122").Value = -4764.8638
$ws.Range("N122").Value = -9451.2001
$ws.Range("H136").Value = 3575.4736
$ws.Range("I136").Value = 2715.6365
$ws.Range("J136").Value = 4757.75
$ws.Range("K136").Value = 8146.9095
$ws.Range("L136").Value = 14273.25
$ws.Range("M136").Value = -5596.9095
$ws.Range("N136").Value = -19373.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 550.36365
$ws.Range("I3").Value = 484.35715
$ws.Range("K3").Value = 484.35715
$ws.Range("M3").Value = -370.35715
$ws.Range("H22").Value = 346.83334
$ws.Range("I22").Value = 360
$ws.Range("J22").Value = 333.66666
$ws.Range("K22").Value = 360
$ws.Range("L22").Value = 333.66666
$ws.Range("M22").Value = -187
$ws.Range("N22").Value = -679.66666
$ws.Range("H80").Value = 1503.5
$ws.Range("J80").Value = 1503.5
$ws.Range("L80").Value = 1503.5
$ws.Range("N80").Value = -3499.5
$ws.Range("H83").Value = 1503.5
$ws.Range("J83").Value = 1503.5
$ws.Range("L83").Value = 7517.5
$ws.Range("N83").Value = -17501.5
$ws.Range("H107").Value = 2044.4375
$ws.Range("I107").Value = 2223.4443
$ws.Range("J107").Value = 1814.2858
$ws.Range("K107").Value = 2223.4443
$ws.Range("L107").Value = 1814.2858
$ws.Range("M107").Value = -303.4443000000001
$ws.Range("N107").Value = -5654.2858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3152.0222
$ws.Range("I31").Value = 1989.8889
$ws.Range("J31").Value = 4895.222
$ws.Range("K31").Value = 1989.8889
$ws.Range("L31").Value = 4895.222
$ws.Range("M31").Value = -1694.8889
$ws.Range("N31").Value = -5485.222
$ws.Range("H34").Value = 3152.0222
$ws.Range("I34").Value = 1989.8889
$ws.Range("J34").Value = 4895.222
$ws.Range("K34").Value = 1989.8889
$ws.Range("L34").Value = 4895.222
$ws.Range("M34").Value = -1787.8889
$ws.Range("N34").Value = -5299.222
$ws.Range("H58").Value = 2091.8064
$ws.Range("I58").Value = 1438.875
$ws.Range("J58").Value = 2788.2666
$ws.Range("K58").Value = 1438.875
$ws.Range("L58").Value = 2788.2666
$ws.Range("M58").Value = -1235.875
$ws.Range("N58").Value = -3194.2666
$ws.Range("H132").Value = 2056.2327
$ws.Range("I132").Value = 1196.7037
$ws.Range("J132").Value = 3506.6875
$ws.Range("K132").Value = 3590.1111
$ws.Range("L132").Value = 10520.0625
$ws.Range("M132").Value = -1060.1111
$ws.Range("N132").Value = -15580.0625
$ws.Range("H136").Value = 2091.8064
$ws.Range("I136").Value = 1438.875
$ws.Range("J136").Value = 2788.2666
$ws.Range("K136").Value = 4316.625
$ws.Range("L136").Value = 8364.799800000001
$ws.Range("M136").Value = -1766.625
$ws.Range("N136").Value = -13464.7998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 53.8
$ws.Range("I12").Value = 17.571428
$ws.Range("K12").Value = 52.71428400000001
$ws.Range("M12").Value = 120.285716
$ws.Range("H56").Value = 2470.3408
$ws.Range("I56").Value = 2470.3408
$ws.Range("K56").Value = 2470.3408
$ws.Range("M56").Value = -1940.3408
$ws.Range("H75").Value = 1105
$ws.Range("I75").Value = 1300
$ws.Range("J75").Value = 1007.5
$ws.Range("K75").Value = 3900
$ws.Range("L75").Value = 3022.5
$ws.Range("M75").Value = -2902
$ws.Range("N75").Value = -5018.5
$ws.Range("H78").Value = 1105
$ws.Range("I78").Value = 1300
$ws.Range("J78").Value = 1007.5
$ws.Range("K78").Value = 11700
$ws.Range("L78").Value = 9067.5
$ws.Range("M78").Value = -6708
$ws.Range("N78").Value = -19051.5
$ws.Range("H86").Value = 170.4
$ws.Range("I86").Value = 134
$ws.Range("J86").Value = 225
$ws.Range("K86").Value = 402
$ws.Range("L86").Value = 675
$ws.Range("M86").Value = 784
$ws.Range("N86").Value = -3047
$ws.Range("H89").Value = 170.4
$ws.Range("I89").Value = 134
$ws.Range("J89").Value = 225
$ws.Range("K89").Value = 1206
$ws.Range("L89").Value = 2025
$ws.Range("M89").Value = 4722
$ws.Range("N89").Value = -13881
$ws.Range("H113").Value = 2632117.2
$ws.Range("I113").Value = 529.98505
$ws.Range("J113").Value = 8929130
$ws.Range("K113").Value = 1589.95515
$ws.Range("L113").Value = 26787390
$ws.Range("M113").Value = 580.04485
$ws.Range("N113").Value = -26791730

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 25293
$ws.Range("J111").Value = 25293
$ws.Range("L111").Value = 25293
$ws.Range("N111").Value = -31427

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H132").Value = 9042.4375
$ws.Range("I132").Value = 2980.2
$ws.Range("J132").Value = 14391.471
$ws.Range("K132").Value = 8940.599999999999
$ws.Range("L132").Value = 43174.413
$ws.Range("M132").Value = -6410.599999999999
$ws.Range("N132").Value = -48234.413
$ws.Range("H136").Value = 4918.6343
$ws.Range("I136").Value = 2846.739
$ws.Range("K136").Value = 8540.217000000001
$ws.Range("M136").Value = -5990.217000000001
$ws.Range("N110").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 10355.5
$ws.Range("I38").Value = 2222
$ws.Range("J38").Value = 13066.667
$ws.Range("K38").Value = 2222
$ws.Range("L38").Value = 13066.667
$ws.Range("M38").Value = -1749
$ws.Range("N38").Value = -14012.667
$ws.Range("H136").Value = 55559228
$ws.Range("I136").Value = 66669812
$ws.Range("J136").Value = 37041590
$ws.Range("K136").Value = 200009436
$ws.Range("L136").Value = 111124770
$ws.Range("M136").Value = -200006886
$ws.Range("N136").Value = -111129870
